# Regen save_data to use K (strikeouts) instead of Strike# for the Trevor Gott
# 2022 game log: recompute column G ("K") for every game row on sheet1.
#
# The workbook has one header row (row 1: date, TB, PC, dS0, dSF, K, IP, I0, IF)
# followed by 51 game rows (rows 2-52). Only column G ("K") values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 2
    12 = 0
    13 = 2
    14 = 0
    15 = 2
    16 = 1
    17 = 2
    18 = 2
    19 = 3
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 0
    33 = 2
    34 = 1
    35 = 0
    36 = 0
    37 = 2
    38 = 2
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    43 = 0
    44 = 1
    45 = 2
    46 = 3
    47 = 0
    48 = 0
    49 = 0
    50 = 2
    51 = 2
    52 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
